# Insert a new data row above the current row 41 (shifting all following
# rows down by one, including the former last row 161 -> 162) and populate
# it with the new weekly record. All "constant" columns (A,B,C,E,F,G,H,N,
# Q,R) match every other row in this sheet; only Fecha (D), Calidad (I)
# and Volumen (J) differ for the new record, while K,L,M,O,P keep the
# values that used to live in (old) row 41.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("41:41").Insert()

$ws.Cells.Item(41, 1).Value = 5
$ws.Cells.Item(41, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(41, 3).Value = "Maule"
$ws.Cells.Item(41, 4).Value = 44414
$ws.Cells.Item(41, 5).Value = 7
$ws.Cells.Item(41, 6).Value = 100112023
$ws.Cells.Item(41, 7).Value = "Brócoli"
$ws.Cells.Item(41, 8).Value = "Sin especificar"
$ws.Cells.Item(41, 9).Value = "Segunda"
$ws.Cells.Item(41, 10).Value = 3000
$ws.Cells.Item(41, 11).Value = 500
$ws.Cells.Item(41, 12).Value = 500
$ws.Cells.Item(41, 13).Value = 500
$ws.Cells.Item(41, 14).Value = "`$/unidad"
$ws.Cells.Item(41, 15).Value = "Región del Maule"
$ws.Cells.Item(41, 16).Value = 500
$ws.Cells.Item(41, 17).Value = 1
$ws.Cells.Item(41, 18).Value = "Hortaliza"
